# Add a new "Spain" test-data sheet, cloned from the "Italy" sheet (same
# layout/styles/merges), placed immediately after "Italy", with the
# market name and NGC code swapped for the Spain-specific values.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Copy "Italy" and place the copy right after "Italy" itself -> Excel
# names it "Italy (2)" and makes it the active sheet, just like using
# the "Move or Copy... > Create a copy" command on the sheet tab.
$italy.Copy($null, $italy)

$spain = $wb.Worksheets.Item("Italy (2)")
$spain.Name = "Spain"

# Set the Spain-specific constant/NGC code first, then the market name,
# so new shared-string entries are minted in that order.
$spain.Range("B4").Value = "NGC-3442/T1599/T1609/T1636"
$spain.Range("B2").Value = "Spain Market"

# Restore Italy's selection to a "whole sheet" selection and make it no
# longer the active tab.
$italy2 = $wb.Worksheets.Item("Italy")
$italy2.Activate()
$italy2.Cells.Select()

# Select A10 on the new Spain sheet and make it the active tab.
$spain2 = $wb.Worksheets.Item("Spain")
$spain2.Activate()
$spain2.Range("A10").Select()
